$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates from the latest cryptos data refresh.
# Numeric-looking Price values are force-written as text (matching the
# original inline-string cell type) by temporarily applying a text
# number format, then restoring the default style so no new formatting
# is left behind on the cell.

$ws.Range("D2").Value = "20.573.13"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").Value = "1.469.56"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9595"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "277.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3560"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3059"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.084"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06622"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.07%  "
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.454"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.164"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9598"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").Value = "1.470.78"
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05958"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.475"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.277"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.69%  "
$ws.Range("D25").Value = "20.576.63"
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.084"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("D29").Value = "1.631.17"
$ws.Range("E29").Value = "  +2.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.951"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.907"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.07946"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7909"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.223"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.91%  "
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05695"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.702"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9604"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.06%  "
$ws.Range("E40").Value = "  +2.20%  "
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1846"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.273"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.514"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5211"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5160"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.795"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06430"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9914"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.29%  "
